$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are written as text, preserving exact formatting
# (e.g. trailing zeros, grouped-dot notation) instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.278.42'
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").Value = '1.841.75'
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '243.49'
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("D6").Value = '0.6872'
$ws.Range("E6").Value = '  -0.99%  '

$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '0.3031'
$ws.Range("E8").Value = '  -0.35%  '

$ws.Range("D9").Value = '0.07524'
$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("D10").Value = '23.28'
$ws.Range("E10").Value = '  +0.15%  '

$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D12").Value = '1.843.03'
$ws.Range("E12").Value = '  +0.61%  '

$ws.Range("D13").Value = '5.089'
$ws.Range("E13").Value = '  -0.01%  '

$ws.Range("D14").Value = '0.6870'
$ws.Range("E14").Value = '  +0.71%  '

$ws.Range("D15").Value = '88.64'
$ws.Range("E15").Value = '  -4.58%  '

$ws.Range("D16").Value = '6.277'
$ws.Range("E16").Value = '  -3.29%  '

$ws.Range("D17").Value = '29.259.61'
$ws.Range("E17").Value = '  +1.22%  '

$ws.Range("D18").Value = '0.000008215'
$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("D19").Value = '2.088.82'
$ws.Range("E19").Value = '  +0.67%  '

$ws.Range("D20").Value = '232.60'
$ws.Range("E20").Value = '  -3.22%  '

$ws.Range("D21").Value = '12.61'
$ws.Range("E21").Value = '  -0.40%  '

$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").Value = '7.475'
$ws.Range("E23").Value = '  +0.38%  '

$ws.Range("D24").Value = '0.9999'
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = '0.1463'
$ws.Range("E25").Value = '  -2.58%  '

$ws.Range("D26").Value = '159.97'
$ws.Range("E26").Value = '  +1.05%  '

$ws.Range("D27").Value = '8.835'
$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("E28").Value = '  -0.28%  '

$ws.Range("D29").Value = '1.521'
$ws.Range("E29").Value = '  -1.15%  '

$ws.Range("D30").Value = '4.278'
$ws.Range("E30").Value = '  +1.26%  '

$ws.Range("D31").Value = '4.162'
$ws.Range("E31").Value = '  +0.57%  '

$ws.Range("D32").Value = '1.212'
$ws.Range("E32").Value = '  +2.24%  '

$ws.Range("D33").Value = '0.05153'
$ws.Range("E33").Value = '  +0.74%  '

$ws.Range("D34").Value = '0.7751'
$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("D35").Value = '1.843'
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("D36").Value = '1.139'
$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").Value = '2.673'
$ws.Range("E37").Value = '  -0.85%  '

$ws.Range("D38").Value = '1.307.18'
$ws.Range("E38").Value = '  +2.66%  '

$ws.Range("D39").Value = '0.01846'
$ws.Range("E39").Value = '  -0.62%  '

$ws.Range("D40").Value = '2.701'
$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").Value = '0.9433'
$ws.Range("E41").Value = '  -1.16%  '

$ws.Range("D42").Value = '105.10'
$ws.Range("E42").Value = '  -1.64%  '

$ws.Range("D43").Value = '5.787'
$ws.Range("E43").Value = '  -5.58%  '

$ws.Range("D44").Value = '0.9997'
$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").Value = '9.695'
$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("D46").Value = '1.985.95'
$ws.Range("E46").Value = '  +0.62%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '65.14'
$ws.Range("E47").Value = '  +2.34%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.5190'
$ws.Range("E48").Value = '  +0.48%  '

$ws.Range("D49").Value = '1.775'
$ws.Range("E49").Value = '  +1.44%  '

$ws.Range("D50").Value = '0.00000000121'
$ws.Range("E50").Value = '  -1.54%  '

$ws.Range("D51").Value = '0.05927'
$ws.Range("E51").Value = '  +0.96%  '
